# Adds a new "as of" forecast column (AJ) to both the "cases" and "deaths"
# sheets, adds a new forecast-target row (48, date 2020-05-21), fills in the
# new nowcast diagonal values, and updates two previously-blank/incorrect
# "Observed" cells (B33 on deaths, B34 on both sheets).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("cases")
$ws2 = $wb.Worksheets.Item("deaths")

foreach ($ws in @($ws1, $ws2)) {

    # --- new header cell AJ1: label is the shared date string "2020-05-07" ---
    # (force text so the engine doesn't silently convert the literal into a
    # date serial number the way typing it into Excel would only do for a
    # *formatted* date cell; every other header cell in the sheet is text)
    $ws.Range("AJ1").NumberFormat = "@"
    $ws.Range("AJ1").Value = "2020-05-07"
    $ws.Range("AJ1").ClearFormats()

    # --- create the (currently empty) AJ cells for existing rows 2-34 ---
    # Copying the style from an already "touched" default-styled neighbour
    # materializes a blank cell without introducing a new cell style.
    $ws.Range("AJ2:AJ34").Style = $ws.Range("AI2").Style

    # --- new row 48 (date 2020-05-21): blank placeholder cells B48:AI48 ---
    $ws.Range("B48:AI48").Style = $ws.Range("AI2").Style

    # --- new row 48's date label, column A (kept as text, like column A
    #     above it, rather than becoming a date serial number) ---
    $ws.Range("A48").NumberFormat = "@"
    $ws.Range("A48").Value = "2020-05-21"
    $ws.Range("A48").ClearFormats()
}

# --- nowcast diagonal values for the new "as of 2020-05-21" column (AJ) ---
# sheet "cases"
$ws1.Range("AJ35").Value = 144602
$ws1.Range("AJ36").Value = 153710
$ws1.Range("AJ37").Value = 161787
$ws1.Range("AJ38").Value = 169926
$ws1.Range("AJ39").Value = 178377
$ws1.Range("AJ40").Value = 185476
$ws1.Range("AJ41").Value = 192765
$ws1.Range("AJ42").Value = 198802
$ws1.Range("AJ43").Value = 204955
$ws1.Range("AJ44").Value = 210362
$ws1.Range("AJ45").Value = 215842
$ws1.Range("AJ46").Value = 221290
$ws1.Range("AJ47").Value = 226070
$ws1.Range("AJ48").Value = 230871

# sheet "deaths"
$ws2.Range("AJ35").Value = 9744
$ws2.Range("AJ36").Value = 10246
$ws2.Range("AJ37").Value = 10812
$ws2.Range("AJ38").Value = 11380
$ws2.Range("AJ39").Value = 12017
$ws2.Range("AJ40").Value = 12535
$ws2.Range("AJ41").Value = 13052
$ws2.Range("AJ42").Value = 13458
$ws2.Range("AJ43").Value = 13914
$ws2.Range("AJ44").Value = 14339
$ws2.Range("AJ45").Value = 14782
$ws2.Range("AJ46").Value = 15158
$ws2.Range("AJ47").Value = 15632
$ws2.Range("AJ48").Value = 15999

# --- Observed-column ("B") corrections ---
# cases: B34 ("Observed" for 2020-05-07) was blank, now populated
$ws1.Range("B34").Value = 135106

# deaths: B33 ("Observed" for 2020-05-06") corrected 8536 -> 8535,
#         and B34 ("Observed" for 2020-05-07") was blank, now populated
$ws2.Range("B33").Value = 8535
$ws2.Range("B34").Value = 9146
